# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4407.5
$ws.Range("I34").Value = 4407.5
$ws.Range("K34").Value = 4407.5
$ws.Range("M34").Value = -4204.5

$ws.Range("H36").Value = 4407.5
$ws.Range("I36").Value = 4407.5
$ws.Range("K36").Value = 4407.5
$ws.Range("M36").Value = -3692.5

$ws.Range("H51").Value = 3098.0952
$ws.Range("I51").Value = 2713.4167
$ws.Range("J51").Value = 3611
$ws.Range("K51").Value = 2713.4167
$ws.Range("L51").Value = 3611
$ws.Range("M51").Value = -2229.4167
$ws.Range("N51").Value = -4579

$ws.Range("H62").Value = 50455.832
$ws.Range("I62").Value = 70017.086
$ws.Range("J62").Value = 11333.333
$ws.Range("K62").Value = 70017.086
$ws.Range("L62").Value = 11333.333
$ws.Range("M62").Value = -69393.086
$ws.Range("N62").Value = -12581.333

$ws.Range("H65").Value = 50455.832
$ws.Range("I65").Value = 70017.086
$ws.Range("J65").Value = 11333.333
$ws.Range("K65").Value = 350085.43
$ws.Range("L65").Value = 56666.665
$ws.Range("M65").Value = -346965.43
$ws.Range("N65").Value = -62906.665

$ws.Range("H69").Value = 200012800
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 200012800
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 600038400
$ws.Range("M69").Value = ""
$ws.Range("N69").Value = -600040148

$ws.Range("H72").Value = 200012800
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 200012800
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 1800115200
$ws.Range("M72").Value = ""
$ws.Range("N72").Value = -1800123936

$ws.Range("H74").Value = 16552.889
$ws.Range("I74").Value = 16552.889
$ws.Range("K74").Value = 16552.889
$ws.Range("M74").Value = -15616.889

$ws.Range("H77").Value = 16552.889
$ws.Range("I77").Value = 16552.889
$ws.Range("K77").Value = 82764.44499999999
$ws.Range("M77").Value = -78084.44499999999

$ws.Range("H116").Value = 5969.8
$ws.Range("I116").Value = 5462.25
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 5462.25
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -2020.25
$ws.Range("N116").Value = -14884

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2140.842
$ws.Range("I74").Value = 2215.5625
$ws.Range("K74").Value = 2215.5625
$ws.Range("M74").Value = -1341.5625

$ws.Range("H77").Value = 2140.842
$ws.Range("I77").Value = 2215.5625
$ws.Range("K77").Value = 11077.8125
$ws.Range("M77").Value = -6709.8125

$ws.Range("H122").Value = 1707.25
$ws.Range("I122").Value = 1369.5
$ws.Range("K122").Value = 4108.5
$ws.Range("M122").Value = -1658.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 31254200
$ws.Range("I86").Value = 4100.9165
$ws.Range("J86").Value = 125004500
$ws.Range("K86").Value = 4100.9165
$ws.Range("L86").Value = 125004500
$ws.Range("M86").Value = -2977.9165
$ws.Range("N86").Value = -125006746

$ws.Range("H89").Value = 31254200
$ws.Range("I89").Value = 4100.9165
$ws.Range("J89").Value = 125004500
$ws.Range("K89").Value = 20504.5825
$ws.Range("L89").Value = 625022500
$ws.Range("M89").Value = -14888.5825
$ws.Range("N89").Value = -625033732

$ws.Range("H105").Value = 3088.8708
$ws.Range("I105").Value = 2567.577
$ws.Range("J105").Value = 5799.6
$ws.Range("K105").Value = 2567.577
$ws.Range("L105").Value = 5799.6
$ws.Range("M105").Value = -820.5770000000002
$ws.Range("N105").Value = -9293.6

$ws.Range("H107").Value = 3279.3547
$ws.Range("I107").Value = 2021.9048
$ws.Range("K107").Value = 2021.9048
$ws.Range("M107").Value = -101.9048

$ws.Range("H134").Value = 2332.087
$ws.Range("I134").Value = 2268.476
$ws.Range("K134").Value = 6805.428
$ws.Range("M134").Value = -4270.428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 224.125
$ws.Range("I7").Value = 24.375
$ws.Range("J7").Value = 423.875
$ws.Range("K7").Value = 24.375
$ws.Range("L7").Value = 423.875
$ws.Range("M7").Value = 88.625
$ws.Range("N7").Value = -649.875

$ws.Range("H31").Value = 6472.385
$ws.Range("I31").Value = 2021.1428
$ws.Range("J31").Value = 11665.5
$ws.Range("K31").Value = 2021.1428
$ws.Range("L31").Value = 11665.5
$ws.Range("M31").Value = -1726.1428
$ws.Range("N31").Value = -12255.5

$ws.Range("H34").Value = 6472.385
$ws.Range("I34").Value = 2021.1428
$ws.Range("J34").Value = 11665.5
$ws.Range("K34").Value = 2021.1428
$ws.Range("L34").Value = 11665.5
$ws.Range("M34").Value = -1819.1428
$ws.Range("N34").Value = -12069.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 33958.234
$ws.Range("I11").Value = 35094.863
$ws.Range("K11").Value = 105284.589
$ws.Range("M11").Value = -105144.589

$ws.Range("H23").Value = 6678.7334
$ws.Range("J23").Value = 5102
$ws.Range("L23").Value = 15306
$ws.Range("N23").Value = -15776

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7339.1333
$ws.Range("I7").Value = 5539.6
$ws.Range("J7").Value = 8238.9
$ws.Range("K7").Value = 5539.6
$ws.Range("L7").Value = 8238.9
$ws.Range("M7").Value = -5427.6
$ws.Range("N7").Value = -8462.9

$ws.Range("H22").Value = 3797.074
$ws.Range("I22").Value = 3806.0908
$ws.Range("J22").Value = 3790.875
$ws.Range("K22").Value = 3806.0908
$ws.Range("L22").Value = 3790.875
$ws.Range("M22").Value = -3511.0908
$ws.Range("N22").Value = -4380.875

$ws.Range("H27").Value = 3797.074
$ws.Range("I27").Value = 3806.0908
$ws.Range("J27").Value = 3790.875
$ws.Range("K27").Value = 3806.0908
$ws.Range("L27").Value = 3790.875
$ws.Range("M27").Value = -3699.0908
$ws.Range("N27").Value = -4004.875

$ws.Range("H68").Value = 5338.04
$ws.Range("I68").Value = 2969.6924
$ws.Range("K68").Value = 2969.6924
$ws.Range("M68").Value = -2220.6924

$ws.Range("H71").Value = 5338.04
$ws.Range("I71").Value = 2969.6924
$ws.Range("K71").Value = 14848.462
$ws.Range("M71").Value = -11104.462

$ws.Range("H122").Value = 4101.4375
$ws.Range("I122").Value = 4602.1816
$ws.Range("J122").Value = 2999.8
$ws.Range("K122").Value = 13806.5448
$ws.Range("L122").Value = 8999.400000000001
$ws.Range("M122").Value = -11356.5448
$ws.Range("N122").Value = -13899.4

$ws.Range("H126").Value = 7339.1333
$ws.Range("I126").Value = 5539.6
$ws.Range("J126").Value = 8238.9
$ws.Range("K126").Value = 16618.8
$ws.Range("L126").Value = 24716.7
$ws.Range("M126").Value = -14148.8
$ws.Range("N126").Value = -29656.7

$ws.Range("H132").Value = 3852.1072
$ws.Range("I132").Value = 3710.55
$ws.Range("K132").Value = 11131.65
$ws.Range("M132").Value = -8601.650000000001

$ws.Range("H136").Value = 13903.134
$ws.Range("I136").Value = 2032
$ws.Range("K136").Value = 6096
$ws.Range("M136").Value = -3546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 40130000
$ws.Range("I5").Value = 500000
$ws.Range("J5").Value = 50037500
$ws.Range("K5").Value = 500000
$ws.Range("L5").Value = 50037500
$ws.Range("M5").Value = -499888
$ws.Range("N5").Value = -50037724

$ws.Range("H39").Value = 23166.666
$ws.Range("I39").Value = 22250
$ws.Range("K39").Value = 22250
$ws.Range("M39").Value = -21837

$ws.Range("H43").Value = 22500
$ws.Range("J43").Value = 22500
$ws.Range("L43").Value = 22500
$ws.Range("N43").Value = -22798

$ws.Range("H45").Value = 9650
$ws.Range("J45").Value = 8580
$ws.Range("L45").Value = 8580
$ws.Range("N45").Value = -9562

$ws.Range("H49").Value = 29255.086
$ws.Range("J49").Value = 28615.385
$ws.Range("L49").Value = 28615.385
$ws.Range("N49").Value = -29075.385

$ws.Range("H122").Value = 4301.478
$ws.Range("J122").Value = 6840.4
$ws.Range("L122").Value = 20521.2
$ws.Range("N122").Value = -25421.2

$ws.Range("H132").Value = 1660.8889
$ws.Range("I132").Value = 1784.0667
$ws.Range("K132").Value = 5352.2001
$ws.Range("M132").Value = -2822.2001

$ws.Range("H136").Value = 985.8095
$ws.Range("I136").Value = 995.85
$ws.Range("K136").Value = 2987.55
$ws.Range("M136").Value = -437.5500000000002
